$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder country names (column A) for the three relocated entries ---
# Marruecos moves up, in front of Kazajistan/Panama (rows 35-37)
$ws.Range("A35").Value = "Marruecos"
$ws.Range("A36").Value = "Kazajistan"
$ws.Range("A37").Value = "Panama"

# Mozambique moves up, in front of Birmania (rows 111-112)
$ws.Range("A111").Value = "Mozambique"
$ws.Range("A112").Value = "Birmania"

# Montserrat moves up, in front of Islas Malvinas (rows 214-215)
$ws.Range("A214").Value = "Montserrat"
$ws.Range("A215").Value = "Islas Malvinas"

# --- Update statistic values (columns B-H) ---
$ws.Range("B4").Value = 7110775
$ws.Range("C4").Value = 12838
$ws.Range("D4").Value = 4366589
$ws.Range("E4").Value = 2538264
$ws.Range("G4").Value = 452
$ws.Range("H4").Value = 205922

$ws.Range("B5").Value = 5700508
$ws.Range("C5").Value = 60012
$ws.Range("D5").Value = 4641811
$ws.Range("E5").Value = 968078
$ws.Range("G5").Value = 598
$ws.Range("H5").Value = 90619

$ws.Range("B6").Value = 4602241
$ws.Range("C6").Value = 6906
$ws.Range("E6").Value = 518204
$ws.Range("G6").Value = 251
$ws.Range("H6").Value = 138410

$ws.Range("B25").Value = 278245
$ws.Range("C25").Value = 1069
$ws.Range("E25").Value = 20840

$ws.Range("B27").Value = 203136
$ws.Range("C27").Value = 9762
$ws.Range("D27").Value = 144686
$ws.Range("E27").Value = 57134

$ws.Range("B31").Value = 129892
$ws.Range("C31").Value = 2249
$ws.Range("E31").Value = 15869
$ws.Range("G31").Value = 45
$ws.Range("H31").Value = 11171

$ws.Range("B35").Value = 107743
$ws.Range("C35").Value = 2397
$ws.Range("D35").Value = 88244
$ws.Range("E35").Value = 17581
$ws.Range("G35").Value = 29
$ws.Range("H35").Value = 1918

$ws.Range("B36").Value = 107450
$ws.Range("C36").Value = 76
$ws.Range("D36").Value = 102064
$ws.Range("E36").Value = 3687
$ws.Range("H36").Value = 1699

$ws.Range("B37").Value = 107284
$ws.Range("D37").Value = 83318
$ws.Range("E37").Value = 21681
$ws.Range("H37").Value = 2285

$ws.Range("B62").Value = 50400
$ws.Range("C62").Value = 186
$ws.Range("D62").Value = 35428
$ws.Range("E62").Value = 13274
$ws.Range("G62").Value = 9
$ws.Range("H62").Value = 1698

$ws.Range("B73").Value = 33675
$ws.Range("C73").Value = 231
$ws.Range("E73").Value = 8517
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = 1794

$ws.Range("B93").Value = 13210
$ws.Range("C93").Value = 57
$ws.Range("E93").Value = 2572

$ws.Range("B100").Value = 9885
$ws.Range("C100").Value = 67
$ws.Range("D100").Value = 8530
$ws.Range("E100").Value = 1321

$ws.Range("B111").Value = 7262
$ws.Range("C111").Value = 148
$ws.Range("D111").Value = 4350
$ws.Range("E111").Value = 2863
$ws.Range("G111").Value = 4
$ws.Range("H111").Value = 49

$ws.Range("B112").Value = 7177
$ws.Range("C112").Value = 434
$ws.Range("D112").Value = 1951
$ws.Range("E112").Value = 5097
$ws.Range("G112").Value = 14
$ws.Range("H112").Value = 129

$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
